$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: AlexNet Implementation -> Silicon Mac TensorFlow install article
$ws.Range("D4").Value = "실리콘(M1, M2) 맥(Mac) OS 에서 텐서플로(TensorFlow) 설치 방법"
$ws.Range("E4").Value = "https://teddylee777.github.io/tensorflow/tensorflow-silicon-install"

# Row 29: title update only (link unchanged)
$ws.Range("D29").Value = "프로메디우스"

# Row 32: Hadoop Ecosystem -> Fact Table / Dimension Table, and link number update
$ws.Range("D32").Value = "Fact Table / Dimension Table"
$ws.Range("E32").Value = "https://dodonam.tistory.com/391"
